$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 110, shifting the existing rows 110-113 down to 112-115
$ws.Rows("110:111").Insert()

# Fill in the two new rows (110, 111) with the "Murcott" entries
$ws.Cells.Item(110, 1).Value = 7
$ws.Cells.Item(110, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(110, 3).Value = "Ñuble"
$ws.Cells.Item(110, 4).Value = 44448
$ws.Cells.Item(110, 5).Value = 16
$ws.Cells.Item(110, 6).Value = "Fruta"
$ws.Cells.Item(110, 7).Value = 100102
$ws.Cells.Item(110, 8).Value = "Cítricos"
$ws.Cells.Item(110, 9).Value = 100102004
$ws.Cells.Item(110, 10).Value = "Mandarina"
$ws.Cells.Item(110, 11).Value = "Murcott"
$ws.Cells.Item(110, 12).Value = "Primera"
$ws.Cells.Item(110, 13).Value = 240
$ws.Cells.Item(110, 14).Value = 6500
$ws.Cells.Item(110, 15).Value = 7000
$ws.Cells.Item(110, 16).Value = 6750
$ws.Cells.Item(110, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(110, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(110, 19).Value = 675
$ws.Cells.Item(110, 20).Value = 10

$ws.Cells.Item(111, 1).Value = 7
$ws.Cells.Item(111, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111, 3).Value = "Ñuble"
$ws.Cells.Item(111, 4).Value = 44448
$ws.Cells.Item(111, 5).Value = 16
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100102
$ws.Cells.Item(111, 8).Value = "Cítricos"
$ws.Cells.Item(111, 9).Value = 100102004
$ws.Cells.Item(111, 10).Value = "Mandarina"
$ws.Cells.Item(111, 11).Value = "Murcott"
$ws.Cells.Item(111, 12).Value = "Segunda"
$ws.Cells.Item(111, 13).Value = 120
$ws.Cells.Item(111, 14).Value = 6000
$ws.Cells.Item(111, 15).Value = 6000
$ws.Cells.Item(111, 16).Value = 6000
$ws.Cells.Item(111, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(111, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(111, 19).Value = 600
$ws.Cells.Item(111, 20).Value = 10
